# Updated cryptos list on Mon Jun 10 17:40:11 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) / Volume(1h) (column E) figures scraped
# for this run. Two pairs of rows (25/26 and 31/32) also changed rank
# order, so their whole row (Coin / Link / Price / Volume) is rewritten.
#
# Price values that look like plain decimals (e.g. "1.00", "6.53") would
# otherwise be auto-parsed as numbers by Excel's normal cell-entry rules,
# so those are entered with a leading apostrophe to force text, exactly
# like typing '1.00 into a cell in the Excel UI. Values already containing
# more than one '.' (e.g. "69.890.98") can't parse as a number and don't
# need the apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function SetRow($r, $coin, $link, $price, $vol) {
    if ($coin -ne "")  { $ws.Cells.Item($r, 2).Value = $coin }
    if ($link -ne "")  { $ws.Cells.Item($r, 3).Value = $link }
    if ($price -ne "") { $ws.Cells.Item($r, 4).Value = $price }
    if ($vol -ne "")   { $ws.Cells.Item($r, 5).Value = $vol }
}

SetRow 2  "" "" "69.890.98"  "  +0.45%  "
SetRow 3  "" "" "3.690.71"   "  -0.11%  "
SetRow 4  "" "" ""           "  -0.08%  "
SetRow 5  "" "" "'647.84"    "  -4.24%  "
SetRow 6  "" "" "'162.02"    "  +0.87%  "
SetRow 7  "" "" ""           "  -0.08%  "
SetRow 8  "" "" ""           "  +1.05%  "
SetRow 9  "" "" ""           "  -0.82%  "
SetRow 10 "" "" "'7.20"      "  +1.40%  "
SetRow 11 "" "" ""           "  +0.68%  "
SetRow 12 "" "" ""           "  -0.28%  "
SetRow 13 "" "" "4.314.86"   "  -0.10%  "
SetRow 15 "" "" "3.675.86"   "  -0.50%  "
SetRow 16 "" "" "69.873.79"  "  +0.39%  "
SetRow 17 "" "" ""           "  +0.63%  "
SetRow 18 "" "" ""           "  -0.10%  "
SetRow 19 "" "" "'6.53"      "  +0.95%  "
SetRow 20 "" "" "'10.38"     "  +5.59%  "
SetRow 21 "" "" "'471.70"    "  +0.38%  "
SetRow 22 "" "" ""           "  +0.41%  "
SetRow 23 "" "" "'80.06"     "  -0.63%  "
SetRow 24 "" "" "3.837.82"   ""

# Rows 25/26 swapped rank order: PEPE now outranks Dai.
SetRow 25 "PEPE" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" "'0.0000127" "  +1.46%  "
SetRow 26 "Dai"  "https://coinranking.com/coin/MoTuySvg7+dai-dai"   "'1.00"      "  +0.01%  "

SetRow 27 "" "" "'10.94" "  +0.55%  "
SetRow 28 "" "" ""       "  +0.45%  "
SetRow 29 "" "" ""       "  -1.64%  "
SetRow 30 "" "" "'1.72"  "  -1.78%  "

# Rows 31/32 swapped rank order: Kaspa now outranks ImmutableX.
SetRow 31 "Kaspa"      "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"      "'0.169" "  +4.19%  "
SetRow 32 "ImmutableX" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" "'2.01"  "  +0.01%  "

SetRow 33 "" "" "'1.00"     "  -0.16%  "
SetRow 34 "" "" "'6.54"     ""
SetRow 35 "" "" "'26.76"    "  -0.65%  "
SetRow 36 "" "" "3.687.19"  "  +0.04%  "
SetRow 37 "" "" ""          "  -0.09%  "
SetRow 38 "" "" ""          "  -0.07%  "
SetRow 39 "" "" "'180.97"   "  +8.43%  "
SetRow 40 "" "" ""          "  -5.06%  "
SetRow 41 "" "" ""          "  -0.07%  "
SetRow 42 "" "" ""          "  -0.19%  "
SetRow 43 "" "" ""          "  +0.32%  "
SetRow 44 "" "" "'0.934"    "  -0.99%  "
SetRow 45 "" "" ""          "  +3.92%  "
SetRow 46 "" "" "'29.40"    "  +4.81%  "
SetRow 47 "" "" "'46.60"    "  -0.87%  "
SetRow 48 "" "" "'0.000275" "  -1.08%  "
SetRow 49 "" "" ""          "  -3.34%  "
SetRow 50 "" "" ""          "  +0.12%  "
SetRow 51 "" "" ""          "  -3.12%  "
